# Apply Trade #15 results to the live_trading_results workbook.

$wb = $excel.ActiveWorkbook

# --- Sheet: Summary ---
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1199.57
$summary.Range("B4").Value = -0.43
$summary.Range("B6").Value = 15
$summary.Range("B8").Value = 8
$summary.Range("B9").Value = 26.67

# --- Sheet: Strategy Status ---
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 99.56999999999999
$status.Range("D4").Value = 15
$status.Range("E4").Value = -0.43
$status.Range("F4").Value = -0.43
$status.Range("G4").Value = 26.67

# --- Sheet: All Trades --- add new trade row 16 ---
$allTrades = $wb.Worksheets.Item("All Trades")
$allTrades.Range("A16").Value = 15
# Force the date-looking text to stay plain text (avoid Excel auto date conversion),
# then clear the forced text format so no stray style gets attached to the cell.
$allTrades.Range("B16").NumberFormat = "@"
$allTrades.Range("B16").Value = "2026-02-17"
$allTrades.Range("B16").ClearFormats()
$allTrades.Range("C16").Value = "08:14:13"
$allTrades.Range("D16").Value = "MarketMaking"
$allTrades.Range("E16").Value = "DOWN"
$allTrades.Range("F16").Value = 0.43
$allTrades.Range("G16").Value = 0.396341
$allTrades.Range("H16").Value = "CLOSED"
$allTrades.Range("I16").Value = -7.8276
$allTrades.Range("J16").Value = -0.03
$allTrades.Range("K16").Value = 99.56999999999999
$allTrades.Range("L16").Value = 0
$allTrades.Range("M16").Value = 0
$allTrades.Range("N16").Value = 0.6
$allTrades.Range("O16").Value = "Normal spread capture: 19600 bps"
$allTrades.Range("P16").Value = "early_exit"
$allTrades.Range("Q16").Value = 0.13

# --- Sheet: MarketMaking --- add new trade row 16 ---
$mm = $wb.Worksheets.Item("MarketMaking")
$mm.Range("A16").Value = 15
$mm.Range("B16").NumberFormat = "@"
$mm.Range("B16").Value = "2026-02-17"
$mm.Range("B16").ClearFormats()
$mm.Range("C16").Value = "08:14:13"
$mm.Range("D16").Value = "MarketMaking"
$mm.Range("E16").Value = "DOWN"
$mm.Range("F16").Value = 0.43
$mm.Range("G16").Value = 0.396341
$mm.Range("H16").Value = "CLOSED"
$mm.Range("I16").Value = -7.8276
$mm.Range("J16").Value = -0.03
$mm.Range("K16").Value = 99.56999999999999
$mm.Range("L16").Value = 0
$mm.Range("M16").Value = 0
$mm.Range("N16").Value = 0.6
$mm.Range("O16").Value = "Normal spread capture: 19600 bps"
$mm.Range("P16").Value = "early_exit"
$mm.Range("Q16").Value = 0.13
